# cryptos.xlsx refresh -- updates Price (D) / Volume(1h) (E) for each coin row,
# and reorders rows 41-42 (EnergySwap now ranks above Stellar).
#
# Set-CellText writes a literal string even when the text parses as a number
# (e.g. "109.82"), mirroring how Excel stores a value typed with a leading
# apostrophe -- otherwise Range.Value would silently coerce it to a Double.
function Set-CellText($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws.Range("D2") '48.688.28'
Set-CellText $ws.Range("E2") '  -2.91%  '

Set-CellText $ws.Range("D3") '2.613.36'
Set-CellText $ws.Range("E3") '  +0.19%  '

Set-CellText $ws.Range("D4") '0.999'
Set-CellText $ws.Range("E4") '  +0.16%  '

Set-CellText $ws.Range("D5") '109.82'
Set-CellText $ws.Range("E5") '  -0.64%  '

Set-CellText $ws.Range("D6") '321.78'
Set-CellText $ws.Range("E6") '  -0.72%  '

Set-CellText $ws.Range("D7") '0.522'
Set-CellText $ws.Range("E7") '  -2.15%  '

Set-CellText $ws.Range("D8") '0.999'
Set-CellText $ws.Range("E8") '  +0.15%  '

Set-CellText $ws.Range("E9") '  -4.14%  '

Set-CellText $ws.Range("D10") '39.34'
Set-CellText $ws.Range("E10") '  -3.64%  '

Set-CellText $ws.Range("D11") '19.71'
Set-CellText $ws.Range("E11") '  -5.21%  '

Set-CellText $ws.Range("E12") '  -2.11%  '

Set-CellText $ws.Range("E13") '  +0.12%  '

Set-CellText $ws.Range("D14") '7.21'
Set-CellText $ws.Range("E14") '  -1.48%  '

Set-CellText $ws.Range("D15") '3.019.36'
Set-CellText $ws.Range("E15") '  +0.13%  '

Set-CellText $ws.Range("D16") '2.607.99'
Set-CellText $ws.Range("E16") '  +1.32%  '

Set-CellText $ws.Range("E17") '  -0.98%  '

Set-CellText $ws.Range("D18") '48.625.36'
Set-CellText $ws.Range("E18") '  -2.65%  '

Set-CellText $ws.Range("E19") '  -4.60%  '

Set-CellText $ws.Range("D20") '12.82'
Set-CellText $ws.Range("E20") '  -4.72%  '

Set-CellText $ws.Range("E21") '  -1.81%  '

Set-CellText $ws.Range("D22") '0.0₃0942'
Set-CellText $ws.Range("E22") '  -1.46%  '

Set-CellText $ws.Range("D23") '269.19'
Set-CellText $ws.Range("E23") '  -5.80%  '

Set-CellText $ws.Range("D24") '68.55'
Set-CellText $ws.Range("E24") '  -6.25%  '

Set-CellText $ws.Range("E25") '  -0.66%  '

Set-CellText $ws.Range("D26") '25.99'
Set-CellText $ws.Range("E26") '  -3.05%  '

Set-CellText $ws.Range("E27") '  +0.03%  '

Set-CellText $ws.Range("D28") '10.00'
Set-CellText $ws.Range("E28") '  -0.05%  '

Set-CellText $ws.Range("E29") '  -0.84%  '

Set-CellText $ws.Range("D30") '34.74'
Set-CellText $ws.Range("E30") '  -3.84%  '

Set-CellText $ws.Range("E31") '  -7.37%  '

Set-CellText $ws.Range("D32") '49.26'
Set-CellText $ws.Range("E32") '  -0.53%  '

Set-CellText $ws.Range("D33") '5.45'
Set-CellText $ws.Range("E33") '  -0.11%  '

Set-CellText $ws.Range("E34") '  -0.07%  '

Set-CellText $ws.Range("E35") '  +0.07%  '

Set-CellText $ws.Range("D36") '18.95'
Set-CellText $ws.Range("E36") '  -4.92%  '

Set-CellText $ws.Range("D37") '4.97'
Set-CellText $ws.Range("E37") '  +4.24%  '

Set-CellText $ws.Range("D38") '2.03'
Set-CellText $ws.Range("E38") '  -2.20%  '

Set-CellText $ws.Range("E39") '  +1.24%  '

Set-CellText $ws.Range("D40") '126.19'
Set-CellText $ws.Range("E40") '  +1.51%  '

Set-CellText $ws.Range("B41") 'EnergySwap'
Set-CellText $ws.Range("C41") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws.Range("D41") '22.45'
Set-CellText $ws.Range("E41") '  -1.35%  '

Set-CellText $ws.Range("B42") 'Stellar'
Set-CellText $ws.Range("C42") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws.Range("D42") '0.111'
Set-CellText $ws.Range("E42") '  -1.81%  '

Set-CellText $ws.Range("E43") '  -4.46%  '

Set-CellText $ws.Range("D44") '0.0316'
Set-CellText $ws.Range("E44") '  +0.44%  '

Set-CellText $ws.Range("D45") '2.059.70'
Set-CellText $ws.Range("E45") '  +0.94%  '

Set-CellText $ws.Range("E46") '  -4.59%  '

Set-CellText $ws.Range("D47") '2.12'
Set-CellText $ws.Range("E47") '  +3.37%  '

Set-CellText $ws.Range("D48") '2.16'
Set-CellText $ws.Range("E48") '  -0.01%  '

Set-CellText $ws.Range("E49") '  -3.58%  '

Set-CellText $ws.Range("D50") '58.45'
Set-CellText $ws.Range("E50") '  +1.40%  '

Set-CellText $ws.Range("D51") '5.15'
Set-CellText $ws.Range("E51") '  -5.05%  '
